# Databook restructure:
#  - Remove the "Program Definitions" and "Metadata" sheets (no longer needed).
#  - On "Parameters", the population label in column A (and, for the very
#    first block, column B as well) is now pulled live from the
#    "Population Definitions" sheet via a formula instead of being a
#    hard-coded label, so edits to population names stay in sync.
#  - On "State Variables", the existing population-label formulas are
#    switched from a relative reference into column B (full name) to an
#    absolute reference into column A (abbreviation).

$wb = $excel.ActiveWorkbook

# --- Remove obsolete sheets -------------------------------------------------
$wb.Worksheets("Program Definitions").Delete()
$wb.Worksheets("Metadata").Delete()

# --- Parameters: make population labels reference Population Definitions ---
$params = $wb.Worksheets("Parameters")

# Row (in Parameters) -> corresponding row in "Population Definitions"
$paramRows = @(2, 3, 4, 5, 8, 9, 10, 11, 14, 15, 16, 17, 20, 21, 22, 23, 26, 27, 28, 29, 32, 33, 34, 35, 38, 39, 40, 41, 44, 45, 46, 47, 50, 51, 52, 53, 56, 57, 58, 59, 62, 63, 64, 65, 68, 69, 70, 71)

foreach ($r in $paramRows) {
    $blockStart = [Math]::Floor(($r - 2) / 6) * 6 + 2
    $offset = $r - $blockStart
    $popRow = 2 + $offset
    $params.Range("A$r").Formula = "='Population Definitions'!A$popRow"
}

# The very first block (rows 2-5) also had column B hard-coded to the same
# full-name text; it now mirrors column A's live formula as well.
foreach ($r in @(2, 3, 4, 5)) {
    $popRow = $r
    $params.Range("B$r").Formula = "='Population Definitions'!A$popRow"
}

# --- State Variables: point the label formula at the abbreviation column ---
$stateVars = $wb.Worksheets("State Variables")

foreach ($r in @(2, 6, 10)) {
    $stateVars.Range("A$r").Formula = "='Population Definitions'!`$A`$2"
}
foreach ($r in @(3, 7, 11)) {
    $stateVars.Range("A$r").Formula = "='Population Definitions'!`$A`$3"
}

# --- Restore a sensible selection on each remaining sheet, ending on
#     "State Variables" so it is left as the active tab. -------------------
$pop = $wb.Worksheets("Population Definitions")
$pop.Activate()
$pop.Range("A2").Select()

$params.Activate()
$params.Range("A2").Select()

$stateVars.Activate()
$stateVars.Range("A15").Select()

Write-Output "done"
